$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# DSM exporter update: the pairwise "parent/child" labels used to be written
# as "child| parent" (row 4 headers) and "parent/ child" (column A labels).
# The exporter now writes them uniformly as "parent.child", and a couple of
# labels that used to carry their own dedicated text ("Component 1.1.1" /
# "Component 1.2") are now abbreviated ("Component 111" / "Component 12").
# Rows 8-11 also now reuse the very same "parent.child" text already used in
# row 4 (E4/F4/G4/H4) instead of keeping their own separate wording.
# ---------------------------------------------------------------------------

# Row 4 - column headers (C..H)
$ws.Range("C4").Value = "Component 1.1.Component 1.1.1"
$ws.Range("D4").Value = "Component 1.Component 1.2"
$ws.Range("E4").Value = "testReqTrace.Component 2"
$ws.Range("F4").Value = "testReqTrace.Component 3"
$ws.Range("G4").Value = "testReqTrace.Component 4"
$ws.Range("H4").Value = "testReqTrace.Component 6"

# Column A - row labels (rows 6..11)
$ws.Range("A6").Value = "Component 1.1.Component 111"
$ws.Range("A7").Value = "Component 1.Component 12"
$ws.Range("A8").Value = "testReqTrace.Component 2"
$ws.Range("A9").Value = "testReqTrace.Component 3"
$ws.Range("A10").Value = "testReqTrace.Component 4"
$ws.Range("A11").Value = "testReqTrace.Component 6"

# ---------------------------------------------------------------------------
# The red/green "change indicator" shading that used to mark the first two
# entries in each of those groups is no longer used - the exporter now
# renders every header/label cell with the same plain (unshaded) look.
# Re-use the already-plain neighbouring cells' formatting so the cells end
# up sharing the very same style instead of getting their own.
# ---------------------------------------------------------------------------

$ws.Range("G4").Copy()
$ws.Range("C4:F4").PasteSpecial(-4122)

$ws.Range("A10").Copy()
$ws.Range("A6:A9").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Column A is a bit narrower now that the labels are shorter.
$ws.Columns.Item(1).ColumnWidth = 30.877604166666668
